$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = -9
$ws.Range("F3").Value = -5
$ws.Range("F5").Value = 5
$ws.Range("F8").Value = 7
$ws.Range("F9").Value = 8
$ws.Range("F10").Value = 1
